$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")

$changes = @(
    @(2, 1, "Última actualización: 14:44:53"),
    @(3, 1, "Total filas: 308"),
    @(13, 3, "16_SANTA ANA"),
    @(14, 3, "17X38_ROMERO"),
    @(52, 1, "06:52:52"),
    @(52, 3, "23_HERNANDEZ"),
    @(52, 4, 70),
    @(53, 1, "07:36:59"),
    @(53, 3, "17_ROMERO"),
    @(53, 4, 26),
    @(81, 1, "07:12:53"),
    @(81, 3, "17X38_ROMERO"),
    @(81, 4, 110),
    @(82, 1, "07:36:59"),
    @(82, 3, "23_HERNANDEZ"),
    @(82, 4, 86),
    @(109, 1, "10:04:17"),
    @(109, 3, "215C_EL PATO"),
    @(109, 4, 0),
    @(110, 1, "08:46:25"),
    @(110, 3, "14_ABASTO"),
    @(110, 4, 78),
    @(158, 1, "10:36:18"),
    @(158, 3, "16_P MOR-SANTA ANA"),
    @(158, 4, 49),
    @(159, 1, "10:04:17"),
    @(159, 3, "11_ETCHEVERRY"),
    @(159, 4, 81),
    @(186, 1, "11:11:31"),
    @(186, 3, "15_ABASTO"),
    @(186, 4, 66),
    @(187, 1, "11:53:59"),
    @(187, 3, "10_OLMOS"),
    @(187, 4, 24),
    @(188, 1, "12:11:45"),
    @(188, 3, "27_EL RETIRO"),
    @(188, 4, 6),
    @(196, 1, "10:48:14"),
    @(196, 3, "11_ETCHEVERRY"),
    @(196, 4, 103),
    @(197, 1, "10:36:18"),
    @(197, 3, "16_P MOR-SANTA ANA"),
    @(197, 4, 115),
    @(206, 1, "10:48:14"),
    @(206, 3, "14_ABASTO"),
    @(206, 4, 115),
    @(207, 1, "10:55:25"),
    @(207, 3, "15X38_ABASTO"),
    @(207, 4, 108),
    @(224, 3, "16_SANTA ANA"),
    @(225, 3, "215_ALUAR"),
    @(234, 1, "12:32:47"),
    @(234, 3, "16_SANTA ANA"),
    @(234, 4, 49),
    @(235, 1, "11:34:25"),
    @(235, 3, "17_ROMERO"),
    @(235, 4, 107),
    @(236, 3, "215A_EL PATO"),
    @(237, 3, "10_OLMOS"),
    @(244, 1, "12:32:47"),
    @(244, 3, "23_HERNANDEZ"),
    @(244, 4, 67),
    @(245, 1, "13:39:24"),
    @(245, 3, "16_SANTA ANA"),
    @(245, 4, 0),
    @(246, 1, "11:53:59"),
    @(246, 3, "17X38_ROMERO"),
    @(246, 4, 106),
    @(273, 1, "13:12:59"),
    @(273, 3, "23_HERNANDEZ"),
    @(273, 4, 79),
    @(274, 3, "15X38_ABASTO"),
    @(275, 3, "15_ABASTO"),
    @(276, 1, "14:31:57"),
    @(276, 3, "16_SANTA ANA"),
    @(276, 4, 0),
    @(283, 1, "14:44:53"),
    @(283, 3, "15X38_ABASTO"),
    @(283, 4, 0),
    @(284, 1, "12:45:57"),
    @(284, 2, "14:44"),
    @(284, 4, 119),
    @(285, 1, "12:53:14"),
    @(285, 2, "14:45"),
    @(285, 3, "215B_EL PATO"),
    @(285, 4, 112),
    @(286, 1, "13:12:59"),
    @(286, 2, "15:00"),
    @(286, 4, 108),
    @(287, 1, "13:39:24"),
    @(287, 2, "15:01"),
    @(287, 3, "81_EL PELIGRO"),
    @(287, 4, 82),
    @(288, 1, "14:10:21"),
    @(288, 2, "15:03"),
    @(288, 3, "23_HERNANDEZ"),
    @(288, 4, 53),
    @(289, 1, "13:51:48"),
    @(289, 2, "15:04"),
    @(289, 4, 73),
    @(290, 1, "13:39:24"),
    @(290, 3, "14_ABASTO"),
    @(290, 4, 94),
    @(291, 1, "14:31:57"),
    @(291, 2, "15:13"),
    @(291, 4, 42),
    @(292, 1, "13:39:24"),
    @(292, 2, "15:14"),
    @(292, 3, "10_OLMOS"),
    @(292, 4, 95),
    @(294, 1, "13:51:48"),
    @(294, 2, "15:24"),
    @(294, 3, "215C_EL PATO"),
    @(294, 4, 93),
    @(295, 3, "16_P MOR-SANTA ANA"),
    @(296, 1, "13:39:24"),
    @(296, 2, "15:25"),
    @(296, 3, "215C_EL PATO"),
    @(296, 4, 106),
    @(297, 1, "14:10:21"),
    @(297, 2, "15:33"),
    @(297, 3, "16_SANTA ANA"),
    @(297, 4, 83),
    @(298, 2, "15:36"),
    @(298, 3, "17X38_ROMERO"),
    @(298, 4, 105),
    @(299, 1, "13:51:48"),
    @(299, 2, "15:44"),
    @(299, 3, "14_ABASTO"),
    @(299, 4, 113),
    @(300, 1, "14:44:53"),
    @(300, 2, "15:44"),
    @(300, 3, "15_ABASTO"),
    @(300, 4, 60),
    @(301, 2, "15:54"),
    @(301, 3, "27_EL RETIRO"),
    @(301, 4, 83),
    @(302, 2, "15:57"),
    @(302, 3, "27_EL RETIRO"),
    @(302, 4, 107),
    @(303, 1, "14:44:53"),
    @(303, 2, "15:57"),
    @(303, 3, "16_SANTA ANA"),
    @(303, 4, 73),
    @(304, 2, "16:00"),
    @(304, 3, "15_ABASTO"),
    @(304, 4, 89),
    @(305, 1, "14:10:21"),
    @(305, 2, "16:01"),
    @(305, 3, "15_ABASTO"),
    @(305, 4, 111),
    @(306, 1, "14:44:53"),
    @(306, 2, "16:10"),
    @(306, 3, "10_OLMOS"),
    @(306, 4, 86),
    @(306, 5, "LP1912"),
    @(307, 1, "14:31:57"),
    @(307, 2, "16:16"),
    @(307, 3, "10_OLMOS"),
    @(307, 4, 105),
    @(307, 5, "LP1912"),
    @(308, 1, "14:31:57"),
    @(308, 2, "16:24"),
    @(308, 3, "11_ETCHEVERRY"),
    @(308, 4, 113),
    @(308, 5, "LP1912"),
    @(309, 1, "14:31:57"),
    @(309, 2, "16:24"),
    @(309, 3, "215_ALUAR"),
    @(309, 4, 113),
    @(309, 5, "LP1912"),
    @(310, 1, "14:44:53"),
    @(310, 2, "16:31"),
    @(310, 3, "16_P MOR-SANTA ANA"),
    @(310, 4, 107),
    @(310, 5, "LP1912"),
    @(311, 1, "14:44:53"),
    @(311, 2, "16:33"),
    @(311, 3, "23_HERNANDEZ"),
    @(311, 4, 109),
    @(311, 5, "LP1912"),
    @(312, 1, "14:44:53"),
    @(312, 2, "16:33"),
    @(312, 3, "225_GOMEZ"),
    @(312, 4, 109),
    @(312, 5, "LP1912"),
    @(313, 1, "14:44:53"),
    @(313, 2, "16:36"),
    @(313, 3, "17X38_ROMERO"),
    @(313, 4, 112),
    @(313, 5, "LP1912"),
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg[0], $chg[1]).Value = $chg[2]
}

$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 14:44:53"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 14:44:53"
